$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 56
$ws.Range("I5").Value = 56
$ws.Range("K5").Value = 56
$ws.Range("M5").Value = 59
$ws.Range("H6").Value = 291.58334
$ws.Range("I6").Value = 291.1111
$ws.Range("K6").Value = 873.3333
$ws.Range("M6").Value = -761.3333
$ws.Range("H18").Value = 469
$ws.Range("J18").Value = 950
$ws.Range("L18").Value = 950
$ws.Range("N18").Value = -1518
$ws.Range("H28").Value = 57149.723
$ws.Range("I28").Value = 64217.75
$ws.Range("K28").Value = 64217.75
$ws.Range("M28").Value = -63732.75
$ws.Range("H40").Value = 6782.357
$ws.Range("I40").Value = 4994.8887
$ws.Range("J40").Value = 9999.799999999999
$ws.Range("K40").Value = 4994.8887
$ws.Range("L40").Value = 9999.799999999999
$ws.Range("M40").Value = -4819.8887
$ws.Range("N40").Value = -10349.8
$ws.Range("H55").Value = 77866.69500000001
$ws.Range("I55").Value = 142973.58
$ws.Range("J55").Value = 1908.6666
$ws.Range("K55").Value = 142973.58
$ws.Range("L55").Value = 1908.6666
$ws.Range("M55").Value = -142759.58
$ws.Range("N55").Value = -2336.6666
$ws.Range("H62").Value = 12506481
$ws.Range("I62").Value = 41670136
$ws.Range("J62").Value = 7771.7144
$ws.Range("K62").Value = 41670136
$ws.Range("L62").Value = 7771.7144
$ws.Range("M62").Value = -41669512
$ws.Range("N62").Value = -9019.714400000001
$ws.Range("H65").Value = 12506481
$ws.Range("I65").Value = 41670136
$ws.Range("J65").Value = 7771.7144
$ws.Range("K65").Value = 208350680
$ws.Range("L65").Value = 38858.572
$ws.Range("M65").Value = -208347560
$ws.Range("N65").Value = -45098.572
$ws.Range("H69").Value = 16250.25
$ws.Range("I69").Value = 5001
$ws.Range("K69").Value = 15003
$ws.Range("M69").Value = -14129
$ws.Range("H72").Value = 16250.25
$ws.Range("I72").Value = 5001
$ws.Range("K72").Value = 45009
$ws.Range("M72").Value = -40641
$ws.Range("H74").Value = 12872.6
$ws.Range("I74").Value = 11006.923
$ws.Range("K74").Value = 11006.923
$ws.Range("M74").Value = -10070.923
$ws.Range("H77").Value = 12872.6
$ws.Range("I77").Value = 11006.923
$ws.Range("K77").Value = 55034.61500000001
$ws.Range("M77").Value = -50354.61500000001
$ws.Range("H86").Value = 2927836
$ws.Range("I86").Value = 3309.375
$ws.Range("J86").Value = 5267457.5
$ws.Range("K86").Value = 3309.375
$ws.Range("L86").Value = 5267457.5
$ws.Range("M86").Value = -2186.375
$ws.Range("N86").Value = -5269703.5
$ws.Range("H88").Value = 2517
$ws.Range("I88").Value = 1998
$ws.Range("J88").Value = 2620.8
$ws.Range("K88").Value = 1998
$ws.Range("L88").Value = 2620.8
$ws.Range("M88").Value = -1592
$ws.Range("N88").Value = -3432.8
$ws.Range("H89").Value = 2927836
$ws.Range("I89").Value = 3309.375
$ws.Range("J89").Value = 5267457.5
$ws.Range("K89").Value = 16546.875
$ws.Range("L89").Value = 26337287.5
$ws.Range("M89").Value = -10930.875
$ws.Range("N89").Value = -26348519.5
$ws.Range("H91").Value = 2517
$ws.Range("I91").Value = 1998
$ws.Range("J91").Value = 2620.8
$ws.Range("K91").Value = 1998
$ws.Range("L91").Value = 2620.8
$ws.Range("M91").Value = -594
$ws.Range("N91").Value = -5428.8
$ws.Range("H98").Value = 1389.9333
$ws.Range("I98").Value = 1026.8462
$ws.Range("K98").Value = 1026.8462
$ws.Range("M98").Value = 471.1538
$ws.Range("H100").Value = 2197.25
$ws.Range("I100").Value = 2130
$ws.Range("J100").Value = 2399
$ws.Range("K100").Value = 2130
$ws.Range("L100").Value = 2399
$ws.Range("M100").Value = -1589
$ws.Range("N100").Value = -3481
$ws.Range("H106").Value = 2322.087
$ws.Range("I106").Value = 2200.5908
$ws.Range("K106").Value = 2200.5908
$ws.Range("M106").Value = -1569.5908
$ws.Range("H107").Value = 63782.75
$ws.Range("I107").Value = 63782.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 63782.75
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -61862.75
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 1389.9333
$ws.Range("I122").Value = 1026.8462
$ws.Range("K122").Value = 3080.5386
$ws.Range("M122").Value = -630.5385999999999
$ws.Range("H137").Value = 8714.817999999999
$ws.Range("I137").Value = 8602.25
$ws.Range("J137").Value = 9015
$ws.Range("K137").Value = 25806.75
$ws.Range("L137").Value = 27045
$ws.Range("M137").Value = -23256.75
$ws.Range("N137").Value = -32145
$ws.Range("H139").Value = 49997.5
$ws.Range("J139").Value = 49997.5
$ws.Range("L139").Value = 49997.5
$ws.Range("N139").Value = -60277.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2944.8164
$ws.Range("I32").Value = 2792.6355
$ws.Range("K32").Value = 2792.6355
$ws.Range("M32").Value = -2505.6355
$ws.Range("H45").Value = 2253.1538
$ws.Range("I45").Value = 1536.375
$ws.Range("K45").Value = 1536.375
$ws.Range("M45").Value = -1159.375
$ws.Range("H61").Value = 5365.9414
$ws.Range("I61").Value = 4894.5
$ws.Range("J61").Value = 6497.4
$ws.Range("K61").Value = 4894.5
$ws.Range("L61").Value = 6497.4
$ws.Range("M61").Value = -4682.5
$ws.Range("N61").Value = -6921.4
$ws.Range("H63").Value = 4296.4116
$ws.Range("I63").Value = 1821.7273
$ws.Range("K63").Value = 1821.7273
$ws.Range("M63").Value = -1135.7273
$ws.Range("H66").Value = 4296.4116
$ws.Range("I66").Value = 1821.7273
$ws.Range("K66").Value = 9108.636500000001
$ws.Range("M66").Value = -5676.636500000001
$ws.Range("H74").Value = 2123.7273
$ws.Range("I74").Value = 1367.8334
$ws.Range("K74").Value = 1367.8334
$ws.Range("M74").Value = -493.8334
$ws.Range("H77").Value = 2123.7273
$ws.Range("I77").Value = 1367.8334
$ws.Range("K77").Value = 6839.166999999999
$ws.Range("M77").Value = -2471.166999999999
$ws.Range("H102").Value = 1965.9231
$ws.Range("I102").Value = 1505.1818
$ws.Range("J102").Value = 4500
$ws.Range("K102").Value = 1505.1818
$ws.Range("L102").Value = 4500
$ws.Range("M102").Value = 116.8181999999999
$ws.Range("N102").Value = -7744
$ws.Range("H122").Value = 6739.5293
$ws.Range("J122").Value = 4418.6665
$ws.Range("L122").Value = 13255.9995
$ws.Range("N122").Value = -18155.9995
$ws.Range("H132").Value = 3412.6667
$ws.Range("I132").Value = 3495.8386
$ws.Range("J132").Value = 3090.375
$ws.Range("K132").Value = 10487.5158
$ws.Range("L132").Value = 9271.125
$ws.Range("M132").Value = -7957.515800000001
$ws.Range("N132").Value = -14331.125
$ws.Range("H133").Value = 61999.8
$ws.Range("J133").Value = 61999.8
$ws.Range("L133").Value = 61999.8
$ws.Range("N133").Value = -67059.8
$ws.Range("H136").Value = 5365.9414
$ws.Range("I136").Value = 4894.5
$ws.Range("J136").Value = 6497.4
$ws.Range("K136").Value = 14683.5
$ws.Range("L136").Value = 19492.2
$ws.Range("M136").Value = -12133.5
$ws.Range("N136").Value = -24592.2
$ws.Range("H138").Value = 62000
$ws.Range("J138").Value = 62000
$ws.Range("L138").Value = 62000
$ws.Range("N138").Value = -72280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2128375.5
$ws.Range("I86").Value = 3403201.2
$ws.Range("J86").Value = 3666
$ws.Range("K86").Value = 3403201.2
$ws.Range("L86").Value = 3666
$ws.Range("M86").Value = -3402078.2
$ws.Range("N86").Value = -5912
$ws.Range("H89").Value = 2128375.5
$ws.Range("I89").Value = 3403201.2
$ws.Range("J89").Value = 3666
$ws.Range("K89").Value = 17016006
$ws.Range("L89").Value = 18330
$ws.Range("M89").Value = -17010390
$ws.Range("N89").Value = -29562
$ws.Range("H94").Value = 1663.6842
$ws.Range("I94").Value = 1970
$ws.Range("K94").Value = 1970
$ws.Range("M94").Value = -1519
$ws.Range("H99").Value = 2191.4707
$ws.Range("I99").Value = 1942.9231
$ws.Range("K99").Value = 1942.9231
$ws.Range("M99").Value = -444.9231
$ws.Range("H107").Value = 717095.0600000001
$ws.Range("I107").Value = 2581.8
$ws.Range("K107").Value = 2581.8
$ws.Range("M107").Value = -661.8000000000002
$ws.Range("H134").Value = 28581.44
$ws.Range("I134").Value = 4159.0835
$ws.Range("J134").Value = 204422.4
$ws.Range("K134").Value = 12477.2505
$ws.Range("L134").Value = 613267.2
$ws.Range("M134").Value = -9942.250499999998
$ws.Range("N134").Value = -618337.2
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 356.5
$ws.Range("I22").Value = 298
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 298
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = 52
$ws.Range("N22").Value = -1700
$ws.Range("H31").Value = 39958.25
$ws.Range("I31").Value = 2693.923
$ws.Range("J31").Value = 72254
$ws.Range("K31").Value = 2693.923
$ws.Range("L31").Value = 72254
$ws.Range("M31").Value = -2398.923
$ws.Range("N31").Value = -72844
$ws.Range("H34").Value = 39958.25
$ws.Range("I34").Value = 2693.923
$ws.Range("J34").Value = 72254
$ws.Range("K34").Value = 2693.923
$ws.Range("L34").Value = 72254
$ws.Range("M34").Value = -2491.923
$ws.Range("N34").Value = -72658
$ws.Range("H59").Value = 28042.875
$ws.Range("J59").Value = 28042.875
$ws.Range("L59").Value = 28042.875
$ws.Range("N59").Value = -30332.875
$ws.Range("H62").Value = 4176
$ws.Range("I62").Value = 2383.3333
$ws.Range("J62").Value = 5968.6665
$ws.Range("K62").Value = 2383.3333
$ws.Range("L62").Value = 5968.6665
$ws.Range("M62").Value = -1759.3333
$ws.Range("N62").Value = -7216.6665
$ws.Range("H65").Value = 4176
$ws.Range("I65").Value = 2383.3333
$ws.Range("J65").Value = 5968.6665
$ws.Range("K65").Value = 11916.6665
$ws.Range("L65").Value = 29843.3325
$ws.Range("M65").Value = -8796.666499999999
$ws.Range("N65").Value = -36083.3325
$ws.Range("H68").Value = 119500
$ws.Range("J68").Value = 119500
$ws.Range("L68").Value = 119500
$ws.Range("N68").Value = -120998
$ws.Range("H71").Value = 119500
$ws.Range("J71").Value = 119500
$ws.Range("L71").Value = 358500
$ws.Range("N71").Value = -365988
$ws.Range("H86").Value = 6400.778
$ws.Range("I86").Value = 5399.2
$ws.Range("K86").Value = 5399.2
$ws.Range("M86").Value = -4276.2
$ws.Range("H89").Value = 6400.778
$ws.Range("I89").Value = 5399.2
$ws.Range("K89").Value = 26996
$ws.Range("M89").Value = -21380
$ws.Range("H122").Value = 4886.5
$ws.Range("J122").Value = 5668
$ws.Range("L122").Value = 17004
$ws.Range("N122").Value = -21904
$ws.Range("H132").Value = 2932.2856
$ws.Range("I132").Value = 2874.75
$ws.Range("J132").Value = 3009
$ws.Range("K132").Value = 8624.25
$ws.Range("L132").Value = 9027
$ws.Range("M132").Value = -6094.25
$ws.Range("N132").Value = -14087

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 800.5
$ws.Range("I7").Value = 701
$ws.Range("K7").Value = 2103
$ws.Range("M7").Value = -1991
$ws.Range("H39").Value = 9716.25
$ws.Range("J39").Value = 14000
$ws.Range("L39").Value = 42000
$ws.Range("N39").Value = -42588
$ws.Range("H92").Value = 1000968.1
$ws.Range("I92").Value = 2000536.4
$ws.Range("J92").Value = 1399.8
$ws.Range("K92").Value = 6001609.199999999
$ws.Range("L92").Value = 4199.4
$ws.Range("M92").Value = -6000361.199999999
$ws.Range("N92").Value = -6695.4
$ws.Range("H120").Value = 9993.5
$ws.Range("I120").Value = 9993.5
$ws.Range("K120").Value = 29980.5
$ws.Range("M120").Value = -25142.5
$ws.Range("H137").Value = 2712.8572
$ws.Range("I137").Value = 2861.818
$ws.Range("J137").Value = 2166.6667
$ws.Range("K137").Value = 8585.454000000002
$ws.Range("L137").Value = 6500.000100000001
$ws.Range("M137").Value = -3485.454000000002
$ws.Range("N137").Value = -16700.0001
$ws.Range("H139").Value = 7259.316
$ws.Range("J139").Value = 7499.643
$ws.Range("L139").Value = 22498.929
$ws.Range("N139").Value = -32778.929

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2146448.2
$ws.Range("I3").Value = 4002424
$ws.Range("K3").Value = 4002424
$ws.Range("M3").Value = -4002308
$ws.Range("H70").Value = 11591
$ws.Range("I70").Value = 9584.375
$ws.Range("J70").Value = 14801.6
$ws.Range("K70").Value = 9584.375
$ws.Range("L70").Value = 14801.6
$ws.Range("M70").Value = -9314.375
$ws.Range("N70").Value = -15341.6
$ws.Range("H73").Value = 11591
$ws.Range("I73").Value = 9584.375
$ws.Range("J73").Value = 14801.6
$ws.Range("K73").Value = 9584.375
$ws.Range("L73").Value = 14801.6
$ws.Range("M73").Value = -8648.375
$ws.Range("N73").Value = -16673.6
$ws.Range("H80").Value = 956186.5
$ws.Range("J80").Value = 1003068.7
$ws.Range("L80").Value = 1003068.7
$ws.Range("N80").Value = -1005064.7
$ws.Range("H83").Value = 956186.5
$ws.Range("J83").Value = 1003068.7
$ws.Range("L83").Value = 5015343.5
$ws.Range("N83").Value = -5025327.5
$ws.Range("H97").Value = 760.9231
$ws.Range("I97").Value = 807.75
$ws.Range("J97").Value = 199
$ws.Range("K97").Value = 807.75
$ws.Range("L97").Value = 199
$ws.Range("M97").Value = -311.75
$ws.Range("N97").Value = -1191
$ws.Range("H102").Value = 2049.2727
$ws.Range("I102").Value = 1167.75
$ws.Range("K102").Value = 1167.75
$ws.Range("M102").Value = 454.25
$ws.Range("H113").Value = 669846.4399999999
$ws.Range("I113").Value = 1668116.9
$ws.Range("J113").Value = 4332.8887
$ws.Range("K113").Value = 1668116.9
$ws.Range("L113").Value = 4332.8887
$ws.Range("M113").Value = -1665946.9
$ws.Range("N113").Value = -8672.8887
$ws.Range("H122").Value = 5312
$ws.Range("I122").Value = 1998.6666
$ws.Range("J122").Value = 7300
$ws.Range("K122").Value = 5995.9998
$ws.Range("L122").Value = 21900
$ws.Range("M122").Value = -3545.9998
$ws.Range("N122").Value = -26800
$ws.Range("H132").Value = 74385.53
$ws.Range("I132").Value = 8521.77
$ws.Range("J132").Value = 502500
$ws.Range("K132").Value = 25565.31
$ws.Range("L132").Value = 1507500
$ws.Range("M132").Value = -23035.31
$ws.Range("N132").Value = -1512560
$ws.Range("H135").Value = 166750000
$ws.Range("J135").Value = 166750000
$ws.Range("L135").Value = 166750000
$ws.Range("N135").Value = -166760140
$ws.Range("H138").Value = 51600
$ws.Range("J138").Value = 51600
$ws.Range("L138").Value = 51600
$ws.Range("N138").Value = -61880

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 200002500
$ws.Range("I2").Value = 400000000
$ws.Range("K2").Value = 400000000
$ws.Range("M2").Value = -399999888
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H35").Value = 10390.25
$ws.Range("I35").Value = 6853.6665
$ws.Range("K35").Value = 6853.6665
$ws.Range("M35").Value = -6517.6665
$ws.Range("H40").Value = 2859.5454
$ws.Range("I40").Value = 2896.7812
$ws.Range("K40").Value = 2896.7812
$ws.Range("M40").Value = -2760.7812
$ws.Range("H46").Value = 2236.5334
$ws.Range("I46").Value = 2377
$ws.Range("J46").Value = 1850.25
$ws.Range("K46").Value = 2377
$ws.Range("L46").Value = 1850.25
$ws.Range("M46").Value = -2189
$ws.Range("N46").Value = -2226.25
$ws.Range("H61").Value = 2967.4119
$ws.Range("I61").Value = 3410.4285
$ws.Range("K61").Value = 3410.4285
$ws.Range("M61").Value = -3208.4285
$ws.Range("H64").Value = 44250
$ws.Range("J64").Value = 44250
$ws.Range("L64").Value = 44250
$ws.Range("N64").Value = -44700
$ws.Range("H67").Value = 44250
$ws.Range("J67").Value = 44250
$ws.Range("L67").Value = 44250
$ws.Range("N67").Value = -45810
$ws.Range("H68").Value = 2779.8
$ws.Range("I68").Value = 2387.25
$ws.Range("K68").Value = 2387.25
$ws.Range("M68").Value = -1638.25
$ws.Range("H71").Value = 2779.8
$ws.Range("I71").Value = 2387.25
$ws.Range("K71").Value = 11936.25
$ws.Range("M71").Value = -8192.25
$ws.Range("H82").Value = 3395.2856
$ws.Range("I82").Value = 3284.4
$ws.Range("K82").Value = 3284.4
$ws.Range("M82").Value = -2923.4
$ws.Range("H85").Value = 3395.2856
$ws.Range("I85").Value = 3284.4
$ws.Range("K85").Value = 3284.4
$ws.Range("M85").Value = -2036.4
$ws.Range("H100").Value = 3207.4167
$ws.Range("I100").Value = 3207.4167
$ws.Range("K100").Value = 3207.4167
$ws.Range("M100").Value = -2666.4167
$ws.Range("H113").Value = 2967.4119
$ws.Range("I113").Value = 3410.4285
$ws.Range("K113").Value = 3410.4285
$ws.Range("M113").Value = -1240.4285
$ws.Range("H122").Value = 3683.3076
$ws.Range("I122").Value = 3262.2727
$ws.Range("K122").Value = 9786.8181
$ws.Range("M122").Value = -7336.8181
$ws.Range("H136").Value = 329097.25
$ws.Range("I136").Value = 560083.6
$ws.Range("K136").Value = 1680250.8
$ws.Range("M136").Value = -1677700.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 1000000000
$ws.Range("J29").Value = 1000000000
$ws.Range("L29").Value = 1000000000
$ws.Range("N29").Value = -1000000580
$ws.Range("H62").Value = 7907.1816
$ws.Range("I62").Value = 7749.75
$ws.Range("K62").Value = 7749.75
$ws.Range("M62").Value = -7125.75
$ws.Range("H63").Value = 120240
$ws.Range("J63").Value = 120240
$ws.Range("L63").Value = 120240
$ws.Range("N63").Value = -121488
$ws.Range("H65").Value = 7907.1816
$ws.Range("I65").Value = 7749.75
$ws.Range("K65").Value = 38748.75
$ws.Range("M65").Value = -35628.75
$ws.Range("H66").Value = 120240
$ws.Range("J66").Value = 120240
$ws.Range("L66").Value = 360720
$ws.Range("N66").Value = -366960
$ws.Range("H81").Value = 4495
$ws.Range("I81").Value = 2494.0908
$ws.Range("J81").Value = 15500
$ws.Range("K81").Value = 4988.1816
$ws.Range("L81").Value = 31000
$ws.Range("M81").Value = -3927.1816
$ws.Range("N81").Value = -33122
$ws.Range("H84").Value = 4495
$ws.Range("I84").Value = 2494.0908
$ws.Range("J84").Value = 15500
$ws.Range("K84").Value = 24940.908
$ws.Range("L84").Value = 155000
$ws.Range("M84").Value = -19636.908
$ws.Range("N84").Value = -165608
$ws.Range("H96").Value = 168467.83
$ws.Range("I96").Value = 504249.5
$ws.Range("K96").Value = 504249.5
$ws.Range("M96").Value = -502876.5
$ws.Range("H100").Value = 683
$ws.Range("I100").Value = 683
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1366
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -825
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 111118450
$ws.Range("I122").Value = 142865500
$ws.Range("J122").Value = 3750.5
$ws.Range("K122").Value = 428596500
$ws.Range("L122").Value = 11251.5
$ws.Range("M122").Value = -428594050
$ws.Range("N122").Value = -16151.5
$ws.Range("H126").Value = 2199
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2199
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 6597
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -11537
$ws.Range("H132").Value = 97338.17999999999
$ws.Range("I132").Value = 6900
$ws.Range("J132").Value = 106382
$ws.Range("K132").Value = 20700
$ws.Range("L132").Value = 319146
$ws.Range("M132").Value = -18170
$ws.Range("N132").Value = -324206
$ws.Range("H133").Value = 47875
$ws.Range("J133").Value = 47875
$ws.Range("L133").Value = 47875
$ws.Range("N133").Value = -57995
